$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: fix separators in "Razon social" name fields (comma -> period) ---
$ws.Range("E82").Value = 'FERNANDEZ. MARIO HUGO'
$ws.Range("E152").Value = 'FERNANDEZ. MARIO HUGO'
$ws.Range("E84").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F84").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E85").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E153").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E108").Value = 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
$ws.Range("F108").Value = 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
$ws.Range("E115").Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
$ws.Range("E177").Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
$ws.Range("E117").Value = 'RICCOTTI. MARIANA EDITH'

# --- Part 2: fix floating point numbers stored as text in "Importe" column ---
# Force text format so numeric-looking strings are not auto-converted to numbers,
# then reset the style back to Normal so no stray formatting is left behind.
$importeRange = $ws.Range("H2:H217")
$importeRange.NumberFormat = "@"

$ws.Range("H2").Value = "3960.11"
$ws.Range("H3").Value = "9106.80"
$ws.Range("H4").Value = "260000.00"
$ws.Range("H5").Value = "436949.15"
$ws.Range("H6").Value = "1168.60"
$ws.Range("H7").Value = "1580.00"
$ws.Range("H8").Value = "1649.54"
$ws.Range("H9").Value = "649.98"
$ws.Range("H10").Value = "5996.99"
$ws.Range("H11").Value = "8605.96"
$ws.Range("H12").Value = "165615.30"
$ws.Range("H13").Value = "10168.00"
$ws.Range("H14").Value = "1900.00"
$ws.Range("H15").Value = "15984.10"
$ws.Range("H16").Value = "8052.00"
$ws.Range("H17").Value = "3985.50"
$ws.Range("H18").Value = "9629.70"
$ws.Range("H19").Value = "970.00"
$ws.Range("H20").Value = "3008.00"
$ws.Range("H21").Value = "800.00"
$ws.Range("H22").Value = "148.00"
$ws.Range("H23").Value = "611.82"
$ws.Range("H24").Value = "822.80"
$ws.Range("H25").Value = "2238.00"
$ws.Range("H26").Value = "886.16"
$ws.Range("H27").Value = "428.00"
$ws.Range("H28").Value = "82.23"
$ws.Range("H29").Value = "23320.00"
$ws.Range("H30").Value = "1232.43"
$ws.Range("H31").Value = "8103.15"
$ws.Range("H32").Value = "1257.11"
$ws.Range("H33").Value = "78.62"
$ws.Range("H34").Value = "2538.50"
$ws.Range("H35").Value = "5849.00"
$ws.Range("H36").Value = "18654.06"
$ws.Range("H37").Value = "77.71"
$ws.Range("H38").Value = "348.00"
$ws.Range("H39").Value = "66.00"
$ws.Range("H40").Value = "55.72"
$ws.Range("H41").Value = "102.75"
$ws.Range("H42").Value = "3064.90"
$ws.Range("H43").Value = "10012.00"
$ws.Range("H44").Value = "30792.29"
$ws.Range("H45").Value = "4719.78"
$ws.Range("H46").Value = "300.00"
$ws.Range("H47").Value = "2135.16"
$ws.Range("H48").Value = "460.00"
$ws.Range("H49").Value = "393.84"
$ws.Range("H50").Value = "27725.42"
$ws.Range("H51").Value = "1110.00"
$ws.Range("H52").Value = "36040.12"
$ws.Range("H53").Value = "536.99"
$ws.Range("H54").Value = "2459.00"
$ws.Range("H55").Value = "180.00"
$ws.Range("H56").Value = "1270.73"
$ws.Range("H57").Value = "1334.40"
$ws.Range("H58").Value = "1008.00"
$ws.Range("H59").Value = "821.10"
$ws.Range("H60").Value = "245.00"
$ws.Range("H61").Value = "5880.00"
$ws.Range("H62").Value = "1446.00"
$ws.Range("H63").Value = "335.00"
$ws.Range("H64").Value = "1300.00"
$ws.Range("H65").Value = "991.98"
$ws.Range("H66").Value = "386.17"
$ws.Range("H67").Value = "289.00"
$ws.Range("H68").Value = "43.92"
$ws.Range("H69").Value = "5382.73"
$ws.Range("H70").Value = "316.11"
$ws.Range("H71").Value = "16667.00"
$ws.Range("H72").Value = "3200.00"
$ws.Range("H73").Value = "2240.00"
$ws.Range("H74").Value = "500.00"
$ws.Range("H75").Value = "1050.00"
$ws.Range("H76").Value = "108972.00"
$ws.Range("H77").Value = "5135.00"
$ws.Range("H78").Value = "4420.00"
$ws.Range("H79").Value = "2420.00"
$ws.Range("H80").Value = "999.00"
$ws.Range("H81").Value = "2697.30"
$ws.Range("H82").Value = "700.00"
$ws.Range("H83").Value = "1139.00"
$ws.Range("H84").Value = "333.59"
$ws.Range("H85").Value = "439.00"
$ws.Range("H86").Value = "40.85"
$ws.Range("H87").Value = "6199.00"
$ws.Range("H88").Value = "1307.67"
$ws.Range("H89").Value = "4500.00"
$ws.Range("H90").Value = "75000.00"
$ws.Range("H91").Value = "7800.00"
$ws.Range("H92").Value = "428940.00"
$ws.Range("H93").Value = "1060.00"
$ws.Range("H94").Value = "706.84"
$ws.Range("H95").Value = "9.41"
$ws.Range("H96").Value = "31.92"
$ws.Range("H97").Value = "635.00"
$ws.Range("H98").Value = "2176.00"
$ws.Range("H99").Value = "378.00"
$ws.Range("H100").Value = "641.79"
$ws.Range("H101").Value = "1129.00"
$ws.Range("H102").Value = "5654.55"
$ws.Range("H103").Value = "100.28"
$ws.Range("H104").Value = "540.00"
$ws.Range("H105").Value = "43.15"
$ws.Range("H106").Value = "9.00"
$ws.Range("H107").Value = "304.00"
$ws.Range("H108").Value = "40.00"
$ws.Range("H109").Value = "129.00"
$ws.Range("H110").Value = "60.00"
$ws.Range("H111").Value = "9701.10"
$ws.Range("H112").Value = "98.00"
$ws.Range("H113").Value = "3249.00"
$ws.Range("H114").Value = "250.00"
$ws.Range("H115").Value = "50.00"
$ws.Range("H116").Value = "60.00"
$ws.Range("H117").Value = "1000.00"
$ws.Range("H118").Value = "200.00"
$ws.Range("H119").Value = "2470.00"
$ws.Range("H120").Value = "2904.00"
$ws.Range("H121").Value = "1620.00"
$ws.Range("H122").Value = "2000.00"
$ws.Range("H123").Value = "4605.76"
$ws.Range("H124").Value = "586.00"
$ws.Range("H125").Value = "2499.36"
$ws.Range("H126").Value = "382.00"
$ws.Range("H127").Value = "2270.30"
$ws.Range("H128").Value = "22.00"
$ws.Range("H129").Value = "45.05"
$ws.Range("H130").Value = "360.00"
$ws.Range("H131").Value = "27.32"
$ws.Range("H132").Value = "131039.20"
$ws.Range("H133").Value = "8160.00"
$ws.Range("H134").Value = "1600.00"
$ws.Range("H135").Value = "1300.00"
$ws.Range("H136").Value = "1000.00"
$ws.Range("H137").Value = "1794.96"
$ws.Range("H138").Value = "384.00"
$ws.Range("H139").Value = "600.00"
$ws.Range("H140").Value = "1000.00"
$ws.Range("H141").Value = "2000.00"
$ws.Range("H142").Value = "4000.00"
$ws.Range("H143").Value = "10588.50"
$ws.Range("H144").Value = "1500.00"
$ws.Range("H145").Value = "950.00"
$ws.Range("H146").Value = "900.00"
$ws.Range("H147").Value = "5000.00"
$ws.Range("H148").Value = "1500.00"
$ws.Range("H149").Value = "200.00"
$ws.Range("H150").Value = "290.00"
$ws.Range("H151").Value = "6480.00"
$ws.Range("H152").Value = "120.00"
$ws.Range("H153").Value = "70.00"
$ws.Range("H154").Value = "395.00"
$ws.Range("H155").Value = "990.00"
$ws.Range("H156").Value = "1800.00"
$ws.Range("H157").Value = "3090.19"
$ws.Range("H158").Value = "3373.24"
$ws.Range("H159").Value = "5910.50"
$ws.Range("H160").Value = "110.00"
$ws.Range("H161").Value = "227.40"
$ws.Range("H162").Value = "4633.00"
$ws.Range("H163").Value = "1942.48"
$ws.Range("H164").Value = "725.00"
$ws.Range("H165").Value = "61.76"
$ws.Range("H166").Value = "3060.00"
$ws.Range("H167").Value = "2886.12"
$ws.Range("H168").Value = "548.40"
$ws.Range("H169").Value = "1395.04"
$ws.Range("H170").Value = "777.03"
$ws.Range("H171").Value = "1007.88"
$ws.Range("H172").Value = "460.66"
$ws.Range("H173").Value = "44.40"
$ws.Range("H174").Value = "2543.00"
$ws.Range("H175").Value = "597.00"
$ws.Range("H176").Value = "15.20"
$ws.Range("H177").Value = "4833.00"
$ws.Range("H178").Value = "129.00"
$ws.Range("H179").Value = "32.00"
$ws.Range("H180").Value = "647.31"
$ws.Range("H181").Value = "79.60"
$ws.Range("H182").Value = "80.00"
$ws.Range("H183").Value = "798.90"
$ws.Range("H184").Value = "966.00"
$ws.Range("H185").Value = "237.50"
$ws.Range("H186").Value = "240.00"
$ws.Range("H187").Value = "352.00"
$ws.Range("H188").Value = "1129.00"
$ws.Range("H189").Value = "114.00"
$ws.Range("H190").Value = "373.91"
$ws.Range("H191").Value = "590.00"
$ws.Range("H192").Value = "995.40"
$ws.Range("H193").Value = "709.40"
$ws.Range("H194").Value = "3958.90"
$ws.Range("H195").Value = "217.80"
$ws.Range("H196").Value = "688864.51"
$ws.Range("H197").Value = "16800.00"
$ws.Range("H198").Value = "1331.00"
$ws.Range("H199").Value = "196300.00"
$ws.Range("H200").Value = "10000.00"
$ws.Range("H201").Value = "52000.00"
$ws.Range("H202").Value = "20000.00"
$ws.Range("H203").Value = "50000.00"
$ws.Range("H204").Value = "223000.00"
$ws.Range("H205").Value = "111500.00"
$ws.Range("H206").Value = "223000.00"
$ws.Range("H207").Value = "19300.00"
$ws.Range("H208").Value = "422794.40"
$ws.Range("H209").Value = "1850.00"
$ws.Range("H210").Value = "188.00"
$ws.Range("H211").Value = "93.98"
$ws.Range("H212").Value = "16000.00"
$ws.Range("H213").Value = "18343.00"
$ws.Range("H214").Value = "8840.00"
$ws.Range("H215").Value = "120000.00"
$ws.Range("H216").Value = "784.00"
$ws.Range("H217").Value = "174.00"

$importeRange.Style = "Normal"
